$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "VSEA" row (current row 10) to hold the new
# "VPOST" / "POSTAGE" office/mode-of-transport pair.
$ws.Range("A10:B10").EntireRow.Insert() | Out-Null

$ws.Range("A10").Value = "VPOST"
$ws.Range("B10").Value = "POSTAGE"

# The newly inserted row inherited the bold header-row style from the row
# above it; reset the font so the new row matches the plain (non-bold)
# look used for its value.
$ws.Range("A10:B10").Font.ThemeColor = 1

# Move the active selection, matching the saved workbook state.
$ws.Range("B14").Select() | Out-Null
